$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Add Sheet2 right after Sheet1 and make it the active sheet
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)

# --- Text (shared-string) cells, written in the exact order needed so that
#     new shared-string table indices line up with the target file ---
$newSheet.Range("C4").Value = "STT"
$newSheet.Range("D4").Value = "Title"
$newSheet.Range("D5").Value = "Nhập môn lập trình"
$newSheet.Range("D6").Value = "Kỹ thuật lập trình"
$newSheet.Range("D7").Value = "Cấu trúc dữ liệu và giải thuật"
$newSheet.Range("D8").Value = "Lập trình hướng đối tượng"
$newSheet.Range("D9").Value = "Mạng máy tính căn bản"
$newSheet.Range("D10").Value = "Hệ điều hành"
$newSheet.Range("D11").Value = "Cơ sở dữ liệu"
$newSheet.Range("D13").Value = "Trí tuệ nhân tạo"
$newSheet.Range("D12").Value = "Lập trình python"
$newSheet.Range("D14").Value = "An toàn thông tin"
$newSheet.Range("D15").Value = "Lập trình web"
$newSheet.Range("D16").Value = "Công nghệ phần mềm"
$newSheet.Range("D19").Value = "Hệ quản trị cơ sở dữ liệu"
$newSheet.Range("D20").Value = "Bảo mật web"
$newSheet.Range("D21").Value = "Thiết kế phần mềm hướng đối tượng"
$newSheet.Range("D22").Value = "Lập trình di động"
$newSheet.Range("C26").Value = "SUM"
$newSheet.Range("D17").Value = "Thương mại điện tử"
$newSheet.Range("D18").Value = "Điện toán đám mây"

# --- numeric cells ---
$newSheet.Range("C5").Value = 1
$newSheet.Range("C6").Value = 2
$newSheet.Range("C7").Value = 3
$newSheet.Range("C8").Value = 4
$newSheet.Range("C9").Value = 5
$newSheet.Range("C10").Value = 6
$newSheet.Range("C11").Value = 7
$newSheet.Range("C12").Value = 8
$newSheet.Range("C13").Value = 9
$newSheet.Range("C14").Value = 10
$newSheet.Range("C15").Value = 11
$newSheet.Range("C16").Value = 12
$newSheet.Range("C17").Value = 13
$newSheet.Range("C18").Value = 14
$newSheet.Range("C19").Value = 13
$newSheet.Range("C20").Value = 14
$newSheet.Range("C21").Value = 15
$newSheet.Range("C22").Value = 16
$newSheet.Range("C23").Value = 17
$newSheet.Range("C24").Value = 18
$newSheet.Range("C25").Value = 19
$newSheet.Range("E6").Value = 1
$newSheet.Range("E11").Value = 1

# --- SUM formula ---
$newSheet.Range("E26").Formula = "=SUM(E5:E25)"

# --- borders around the SUM row (C26:E26) ---
$newSheet.Range("C26").Borders.Item(8).LineStyle = 1
$newSheet.Range("C26").Borders.Item(9).LineStyle = 1
$newSheet.Range("C26").Borders.Item(7).LineStyle = 1
$newSheet.Range("E26").Borders.Item(8).LineStyle = 1
$newSheet.Range("E26").Borders.Item(9).LineStyle = 1
$newSheet.Range("E26").Borders.Item(10).LineStyle = 1
$newSheet.Range("D26").Borders.Item(8).LineStyle = 1
$newSheet.Range("D26").Borders.Item(9).LineStyle = 1

# --- column width for column D ---
$newSheet.Columns.Item(4).ColumnWidth = 23.7142857

# --- page setup ---
$newSheet.PageSetup.Orientation = 1

# --- view / selection ---
$newSheet.Range("H20").Select() | Out-Null
